# Fix typo in resume
#
# 1. Remove the stray "_GoBack" bookmark that sits after the
#    "...prevent eavesdropping." paragraph.
# 2. Split "Drove development of Microsoft HoloLens app " into
#    "Drove development of" + " a" + " Microsoft HoloLens app " so the
#    bullet reads "Drove development of a Microsoft HoloLens app...".
# 3. Word moves the (hidden) "_GoBack" bookmark to the most recently
#    edited spot whenever the document is saved, so re-create it as a
#    zero-length bookmark right after "...using this microphone."

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -----------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. "Drove development of Microsoft HoloLens app " -> 3 runs ----
$holo = $d.Paragraphs(19)
$holoText = $holo.Range.Text
$holoStart = $holo.Range.Start

# Insert " a" right before " Microsoft" (turns "...development of
# Microsoft..." into "...development of a Microsoft...").
$msOffset = $holoText.IndexOf(" Microsoft")
$insertPoint = $holoStart + $msOffset
$ip = $d.Range($insertPoint, $insertPoint)
$ip.InsertBefore(" a")

# Force the freshly inserted " a" to stay in its own run instead of
# being re-merged into its neighbours by toggling (and clearing) a
# character property over exactly that span.
$aSeg = $d.Range($insertPoint, $insertPoint + 2)
$aSeg.Font.Bold = $true
$aSeg.Font.Bold = $false

# Likewise keep " Microsoft HoloLens app " as its own run, distinct
# from the following "allowing the hearing-impaired..." run.
$holoText2 = $holo.Range.Text
$msStart = $holoStart + $holoText2.IndexOf(" Microsoft")
$msEnd = $holoStart + $holoText2.IndexOf("allowing")
$msSeg = $d.Range($msStart, $msEnd)
$msSeg.Font.Bold = $true
$msSeg.Font.Bold = $false

# --- 3. Re-add _GoBack at the end of the microphone bullet ----------
$mic = $d.Paragraphs(21)
$micEnd = $mic.Range.End - 1

# Adding a zero-length bookmark exactly at a paragraph's last valid
# offset drops it at the wrong spot, so nudge a placeholder character
# in first, anchor the bookmark next to it, then remove the
# placeholder again -- the bookmark (an independent, zero-length
# marker) stays put in the correct location.
$placeholder = $d.Range($micEnd, $micEnd)
$placeholder.InsertAfter("X")
$bmRange = $d.Range($micEnd, $micEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($micEnd, $micEnd + 1).Delete()
